$wb = $excel.ActiveWorkbook

# Source sheet to base the new sheet's look on
$src = $wb.Worksheets.Item(1)

# Add the new sheet after the last existing sheet (EventsChart) and rename it
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$new.Name = "HailStormCalc"

# Copy cell formatting (styles only, no values/formulas) from RidersCalc's
# data block so the new sheet inherits the same look (header row style,
# first/default numeric column styles, etc.)
$src.Range("A3:D28").Copy()
$new.Range("A3").PasteSpecial(-4122)

# Column A beyond the data (rows 16-28) should have no cell at all, matching
# the source workbook layout where mileage stops part-way down the sheet.
$new.Range("A16:A28").Clear()

# Header row (code comment) text, with its own font/alignment style.
# Prime the cell with the workbook's existing "mono font" style (copied from
# column A of RidersCalc) so that changing the alignment reuses the existing
# font entry instead of fabricating a brand new (mismatched) one.
$src.Range("A5").Copy()
$new.Range("A1").PasteSpecial(-4122)
$new.Range("A1").Value = "IF RND(-1)*10>9-((M/100-15)**2+72)/((M/100-15)**2+12)"
$new.Range("A1").HorizontalAlignment = -4131

# Column headers (row 3)
$new.Range("A3").Value = "Mileage"
$new.Range("B3").Value = "Numerator"
$new.Range("C3").Value = "Denominator"
$new.Range("D3").Value = "Quotient"

# Data rows 4-15: mileage 950..2050 in steps of 100
$new.Range("A4").Value = 950
for ($r = 5; $r -le 15; $r++) {
    $prev = $r - 1
    $new.Range("A$r").Formula = "=A$prev+100"
}
for ($r = 4; $r -le 15; $r++) {
    $new.Range("B$r").Formula = "=(A$r/100 - 15) ^ 2 + 72"
    $new.Range("C$r").Formula = "=(A$r/100 -15) ^ 2 + 12"
    $new.Range("D$r").Formula = "=B$r/C$r"
}

# Column widths, matching RidersCalc (column A keeps the sheet default width,
# columns B:D are widened, same as the source sheet)
$new.Range("B1:D1").EntireColumn.ColumnWidth = 18.616666666666667

# Selection / view state for the new sheet
$new.Range("A1").Select()
